# Reorders the fantasy-roster rows so each player lines up with their
# correct position/team again (rows 4-16 of the "Oyuncu Adı / Pozisyon /
# Takım" table get shuffled; rows 2,3,17,18,19 stay put).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Dejounte Murray",    "PG,SG",   "New Orleans Pelicans"),
    @("Russell Westbrook",  "PG,SG",   "Denver Nuggets"),
    @("Keon Johnson",       "PG,SG",   "Brooklyn Nets"),
    @("Chris Paul",         "PG",      "San Antonio Spurs"),
    @("Pascal Siakam",      "SF,PF,C", "Indiana Pacers"),
    @("Khris Middleton",    "SF",      "Milwaukee Bucks"),
    @("Deni Avdija",        "SF,PF",   "Portland Trail Blazers"),
    @("Nikola Jokic",       "C",       "Denver Nuggets"),
    @("Jakob Poeltl",       "C",       "Toronto Raptors"),
    @("Jonas Valanciunas",  "C",       "Washington Wizards"),
    @("Rudy Gobert",        "C",       "Minnesota Timberwolves"),
    @("Jaylen Brown",       "SG,SF",   "Boston Celtics"),
    @("Tyus Jones",         "PG",      "Phoenix Suns"),
    @("Jalen Green",        "PG,SG",   "Houston Rockets"),
    @("Draymond Green",     "PF,C",    "Golden State Warriors"),
    @("Paolo Banchero",     "SF,PF",   "Orlando Magic"),
    @("Chet Holmgren",      "PF,C",    "Oklahoma City Thunder"),
    @("Jalen Suggs",        "PG,SG",   "Orlando Magic")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row++
}
